$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.917.53'
$ws.Range('E2').Value = '  +1.58%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.407.29'
$ws.Range('E3').Value = '  +1.77%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.93'
$ws.Range('E5').Value = '  +1.32%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.92'
$ws.Range('E6').Value = '  +3.34%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.528'
$ws.Range('E8').Value = '  +0.74%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.401.10'
$ws.Range('E9').Value = '  +1.51%  '

$ws.Range('E10').Value = '  +1.49%  '

$ws.Range('E11').Value = '  -0.89%  '

$ws.Range('E12').Value = '  +0.93%  '

$ws.Range('E13').Value = '  +1.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.88'
$ws.Range('E14').Value = '  +3.81%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000173'
$ws.Range('E15').Value = '  +5.20%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.839.53'
$ws.Range('E16').Value = '  +2.29%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.906.43'
$ws.Range('E17').Value = '  +1.67%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.405.10'
$ws.Range('E18').Value = '  +1.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.06'
$ws.Range('E19').Value = '  +2.98%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.96'
$ws.Range('E20').Value = '  +1.02%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.16'
$ws.Range('E21').Value = '  +1.18%  '

$ws.Range('E22').Value = '  +0.87%  '

$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.98'
$ws.Range('E24').Value = '  +1.55%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.70'
$ws.Range('E25').Value = '  +2.80%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.01'
$ws.Range('E26').Value = '  +9.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '572.08'
$ws.Range('E27').Value = '  +14.51%  '

$ws.Range('E28').Value = '  +0.35%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.526.56'
$ws.Range('E29').Value = '  +2.10%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.22'
$ws.Range('E30').Value = '  +1.81%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0919'
$ws.Range('E31').Value = '  +5.03%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.43'
$ws.Range('E32').Value = '  +4.42%  '

$ws.Range('E33').Value = '  -0.25%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  +3.31%  '

$ws.Range('E35').Value = '  +3.32%  '

$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.61'
$ws.Range('E37').Value = '  +5.87%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.73'
$ws.Range('E38').Value = '  +1.48%  '

$ws.Range('E39').Value = '  +1.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '150.80'
$ws.Range('E40').Value = '  +3.83%  '

$ws.Range('E41').Value = '  +0.43%  '

$ws.Range('E42').Value = '  -2.53%  '

$ws.Range('E43').Value = '  +0.04%  '

$ws.Range('E44').Value = '  +12.66%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '149.52'
$ws.Range('E45').Value = '  +2.22%  '

$ws.Range('E46').Value = '  +1.13%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0535'
$ws.Range('E47').Value = '  +3.17%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.01'
$ws.Range('E48').Value = '  +4.49%  '

$ws.Range('E49').Value = '  +2.26%  '

$ws.Range('E50').Value = '  +1.86%  '

$ws.Range('E51').Value = '  +2.22%  '
